# modify part1's rf: adjust rf to the same period of return
# (previously annual rf was used in all circumstances, which was not
# correct). This updates the GRS test output values (A:E, rows 2-6)
# on the active sheet to reflect the corrected risk-free rate period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tiny p-values below are expressed via mantissa * 10^exponent since
# literal scientific notation (e.g. 1e-16) is not parsed here.
$tinyB = 1.110223024625157 * [Math]::Pow(10, -16)
$tinyB6 = 11.65734175856414 * [Math]::Pow(10, -15)

# Row 2
$ws.Range("A2").Value = 16.61753674131967
$ws.Range("B2").Value = $tinyB
$ws.Range("C2").Value = 0.01008485813524562
$ws.Range("D2").Value = 0.6111179909053543
$ws.Range("E2").Value = 0.3734651988081968

# Row 3
$ws.Range("A3").Value = 12.0063920414376
$ws.Range("B3").Value = $tinyB
$ws.Range("C3").Value = 0.007484340220086571
$ws.Range("D3").Value = 0.4535329002364838
$ws.Range("E3").Value = 0.2056920915969164

# Row 4
$ws.Range("A4").Value = 14.61654013638949
$ws.Range("B4").Value = $tinyB
$ws.Range("C4").Value = 0.01136024991398938
$ws.Range("D4").Value = 0.6884036453975196
$ws.Range("E4").Value = 0.4738995789965941

# Row 5
$ws.Range("A5").Value = 13.95865213269506
$ws.Range("B5").Value = $tinyB
$ws.Range("C5").Value = 0.01093531990953677
$ws.Range("D5").Value = 0.6626539157420355
$ws.Range("E5").Value = 0.4391102120482526

# Row 6
$ws.Range("A6").Value = 10.23090431328895
$ws.Range("B6").Value = $tinyB6
$ws.Range("C6").Value = 0.004541505896120918
$ws.Range("D6").Value = 0.2752042638282132
$ws.Range("E6").Value = 0.07573738682922874
